$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 4781.8423
$ws.Range("I9").Value = 6042.3335
$ws.Range("J9").Value = 55
$ws.Range("K9").Value = 6042.3335
$ws.Range("L9").Value = 55
$ws.Range("M9").Value = -5873.3335
$ws.Range("N9").Value = -393
$ws.Range("H43").Value = 1069
$ws.Range("I43").Value = 973.25
$ws.Range("K43").Value = 973.25
$ws.Range("M43").Value = -904.25
$ws.Range("H88").Value = 1104.625
$ws.Range("J88").Value = 1126.8182
$ws.Range("L88").Value = 1126.8182
$ws.Range("N88").Value = -1938.8182
$ws.Range("H91").Value = 1104.625
$ws.Range("J91").Value = 1126.8182
$ws.Range("L91").Value = 1126.8182
$ws.Range("N91").Value = -3934.8182
$ws.Range("H107").Value = 1032.1666
$ws.Range("I107").Value = 1008.8571
$ws.Range("J107").Value = 1195.3334
$ws.Range("K107").Value = 1008.8571
$ws.Range("L107").Value = 1195.3334
$ws.Range("M107").Value = 911.1429000000001
$ws.Range("N107").Value = -5035.3334
$ws.Range("H112").Value = 113965.11
$ws.Range("I112").Value = 2646.75
$ws.Range("J112").Value = 203019.8
$ws.Range("K112").Value = 7940.25
$ws.Range("L112").Value = 609059.3999999999
$ws.Range("M112").Value = -6832.25
$ws.Range("N112").Value = -611275.3999999999
$ws.Range("H125").Value = 8048.421
$ws.Range("J125").Value = 7214.143
$ws.Range("L125").Value = 64927.287
$ws.Range("N125").Value = -69847.287
$ws.Range("H132").Value = 1515.7142
$ws.Range("I132").Value = 1483.4426
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 4450.3278
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -1920.3278
$ws.Range("N132").Value = -12560
$ws.Range("H137").Value = 6133.1665
$ws.Range("I137").Value = 6975
$ws.Range("J137").Value = 4449.5
$ws.Range("K137").Value = 20925
$ws.Range("L137").Value = 13348.5
$ws.Range("M137").Value = -18375
$ws.Range("N137").Value = -18448.5
$ws.Range("H138").Value = 7149806
$ws.Range("J138").Value = 7469757
$ws.Range("L138").Value = 22409271
$ws.Range("N138").Value = -22419551

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 10823.818
$ws.Range("J31").Value = 29758
$ws.Range("L31").Value = 29758
$ws.Range("N31").Value = -30346
$ws.Range("H32").Value = 22807.896
$ws.Range("I32").Value = 21081.412
$ws.Range("K32").Value = 21081.412
$ws.Range("M32").Value = -20794.412
$ws.Range("H74").Value = 26402.846
$ws.Range("J74").Value = 54000
$ws.Range("L74").Value = 54000
$ws.Range("N74").Value = -55748
$ws.Range("H77").Value = 26402.846
$ws.Range("J77").Value = 54000
$ws.Range("L77").Value = 270000
$ws.Range("N77").Value = -278736
$ws.Range("H94").Value = 29469.666
$ws.Range("J94").Value = 29469.666
$ws.Range("L94").Value = 29469.666
$ws.Range("N94").Value = -31271.666
$ws.Range("H98").Value = 44088.75
$ws.Range("J98").Value = 44088.75
$ws.Range("L98").Value = 44088.75
$ws.Range("N98").Value = -50078.75
$ws.Range("H101").Value = 50000
$ws.Range("J101").Value = 50000
$ws.Range("L101").Value = 50000
$ws.Range("N101").Value = -56490
$ws.Range("H125").Value = 29999
$ws.Range("J125").Value = 29999
$ws.Range("L125").Value = 29999
$ws.Range("N125").Value = -39839
$ws.Range("H132").Value = 3544.4443
$ws.Range("I132").Value = 3780
$ws.Range("K132").Value = 11340
$ws.Range("M132").Value = -8810

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 591.125
$ws.Range("I80").Value = 420.42856
$ws.Range("J80").Value = 723.8889
$ws.Range("K80").Value = 420.42856
$ws.Range("L80").Value = 723.8889
$ws.Range("M80").Value = 577.5714399999999
$ws.Range("N80").Value = -2719.8889
$ws.Range("H83").Value = 591.125
$ws.Range("I83").Value = 420.42856
$ws.Range("J83").Value = 723.8889
$ws.Range("K83").Value = 2102.1428
$ws.Range("L83").Value = 3619.4445
$ws.Range("M83").Value = 2889.8572
$ws.Range("N83").Value = -13603.4445
$ws.Range("H86").Value = 3173.5334
$ws.Range("I86").Value = 3008.2727
$ws.Range("J86").Value = 3628
$ws.Range("K86").Value = 3008.2727
$ws.Range("L86").Value = 3628
$ws.Range("M86").Value = -1885.2727
$ws.Range("N86").Value = -5874
$ws.Range("H89").Value = 3173.5334
$ws.Range("I89").Value = 3008.2727
$ws.Range("J89").Value = 3628
$ws.Range("K89").Value = 15041.3635
$ws.Range("L89").Value = 18140
$ws.Range("M89").Value = -9425.363499999999
$ws.Range("N89").Value = -29372
$ws.Range("H107").Value = 1036.8096
$ws.Range("I107").Value = 1054.421
$ws.Range("K107").Value = 1054.421
$ws.Range("M107").Value = 865.579

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 39999
$ws.Range("J18").Value = 39999
$ws.Range("L18").Value = 39999
$ws.Range("N18").Value = -40459
$ws.Range("H31").Value = 3797.4375
$ws.Range("I31").Value = 3583.9333
$ws.Range("K31").Value = 3583.9333
$ws.Range("M31").Value = -3288.9333
$ws.Range("H34").Value = 3797.4375
$ws.Range("I34").Value = 3583.9333
$ws.Range("K34").Value = 3583.9333
$ws.Range("M34").Value = -3381.9333
$ws.Range("H58").Value = 5003.6665
$ws.Range("I58").Value = 5003.6665
$ws.Range("K58").Value = 5003.6665
$ws.Range("M58").Value = -4800.6665
$ws.Range("H64").Value = 11750
$ws.Range("J64").Value = 12500
$ws.Range("L64").Value = 12500
$ws.Range("N64").Value = -12996
$ws.Range("H67").Value = 11750
$ws.Range("J67").Value = 12500
$ws.Range("L67").Value = 12500
$ws.Range("N67").Value = -14216
$ws.Range("H99").Value = 2708.1667
$ws.Range("I99").Value = 2250
$ws.Range("K99").Value = 2250
$ws.Range("M99").Value = -752
$ws.Range("H122").Value = 3272.7334
$ws.Range("I122").Value = 1514.4445
$ws.Range("J122").Value = 5910.1665
$ws.Range("K122").Value = 4543.333500000001
$ws.Range("L122").Value = 17730.4995
$ws.Range("M122").Value = -2093.333500000001
$ws.Range("N122").Value = -22630.4995
$ws.Range("H126").Value = 2708.1667
$ws.Range("I126").Value = 2250
$ws.Range("K126").Value = 6750
$ws.Range("M126").Value = -4280
$ws.Range("H134").Value = 13419.272
$ws.Range("I134").Value = 7813.1763
$ws.Range("K134").Value = 23439.5289
$ws.Range("M134").Value = -20904.5289
$ws.Range("H136").Value = 5003.6665
$ws.Range("I136").Value = 5003.6665
$ws.Range("K136").Value = 15010.9995
$ws.Range("M136").Value = -12460.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 112889.445
$ws.Range("I14").Value = 112889.445
$ws.Range("K14").Value = 338668.335
$ws.Range("M14").Value = -338495.335
$ws.Range("H129").Value = 14143418
$ws.Range("I129").Value = 14143418
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 42430254
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = -42425254
$ws.Range("N129").ClearContents()
$ws.Range("H131").Value = 3787.0435
$ws.Range("J131").Value = 5000.0625
$ws.Range("L131").Value = 15000.1875
$ws.Range("N131").Value = -25080.1875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 146.08
$ws.Range("I2").Value = 176.2
$ws.Range("K2").Value = 176.2
$ws.Range("M2").Value = -63.19999999999999
$ws.Range("H80").Value = 6874.8335
$ws.Range("I80").Value = 4630
$ws.Range("J80").Value = 7997.25
$ws.Range("K80").Value = 4630
$ws.Range("L80").Value = 7997.25
$ws.Range("M80").Value = -3632
$ws.Range("N80").Value = -9993.25
$ws.Range("H83").Value = 6874.8335
$ws.Range("I83").Value = 4630
$ws.Range("J83").Value = 7997.25
$ws.Range("K83").Value = 23150
$ws.Range("L83").Value = 39986.25
$ws.Range("M83").Value = -18158
$ws.Range("N83").Value = -49970.25
$ws.Range("H126").Value = 14516.8
$ws.Range("J126").Value = 4699.4
$ws.Range("L126").Value = 14098.2
$ws.Range("N126").Value = -19038.2
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H135").Value = 49750
$ws.Range("J135").Value = 49750
$ws.Range("L135").Value = 49750
$ws.Range("N135").Value = -59890
$ws.Range("H136").Value = 4451.3423
$ws.Range("I136").Value = 4365.3057
$ws.Range("J136").Value = 6000
$ws.Range("K136").Value = 13095.9171
$ws.Range("L136").Value = 18000
$ws.Range("M136").Value = -10545.9171
$ws.Range("N136").Value = -23100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 13846.571
$ws.Range("I63").Value = 11335.2
$ws.Range("J63").Value = 20125
$ws.Range("K63").Value = 11335.2
$ws.Range("L63").Value = 20125
$ws.Range("M63").Value = -10711.2
$ws.Range("N63").Value = -21373
$ws.Range("H66").Value = 13846.571
$ws.Range("I66").Value = 11335.2
$ws.Range("J66").Value = 20125
$ws.Range("K66").Value = 34005.60000000001
$ws.Range("L66").Value = 60375
$ws.Range("M66").Value = -30885.60000000001
$ws.Range("N66").Value = -66615
$ws.Range("H132").Value = 4081.6072
$ws.Range("I132").Value = 3904.9443
$ws.Range("J132").Value = 4399.6
$ws.Range("K132").Value = 11714.8329
$ws.Range("L132").Value = 13198.8
$ws.Range("M132").Value = -9184.832900000001
$ws.Range("N132").Value = -18258.8
$ws.Range("H136").Value = 1851.0286
$ws.Range("I136").Value = 1827.2759
$ws.Range("J136").Value = 1965.8334
$ws.Range("K136").Value = 5481.8277
$ws.Range("L136").Value = 5897.5002
$ws.Range("M136").Value = -2931.8277
$ws.Range("N136").Value = -10997.5002
